$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 116
$ws.Range("H116").Value = 2400.8572
$ws.Range("I116").Value = 1638.6364
$ws.Range("J116").Value = 3690.7693
$ws.Range("K116").Value = 1638.6364
$ws.Range("L116").Value = 3690.7693
$ws.Range("M116").Value = 1803.3636
$ws.Range("N116").Value = -10574.7693

# Row 137
$ws.Range("H137").Value = 6170.6
$ws.Range("I137").Value = 7400.421
$ws.Range("J137").Value = 5416.839
$ws.Range("K137").Value = 22201.263
$ws.Range("L137").Value = 16250.517
$ws.Range("M137").Value = -19651.263
$ws.Range("N137").Value = -21350.517

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 52633700
$ws.Range("I45").Value = 90910840
$ws.Range("J45").Value = 2633
$ws.Range("K45").Value = 90910840
$ws.Range("L45").Value = 2633
$ws.Range("M45").Value = -90910463
$ws.Range("N45").Value = -3387

# Row 61
$ws.Range("H61").Value = 1503.7368
$ws.Range("I61").Value = 1017.25
$ws.Range("J61").Value = 2126.44
$ws.Range("K61").Value = 1017.25
$ws.Range("L61").Value = 2126.44
$ws.Range("M61").Value = -805.25
$ws.Range("N61").Value = -2550.44

# Row 110
$ws.Range("H110").Value = 2031.2
$ws.Range("I110").Value = 2090.2222
$ws.Range("J110").Value = 1500
$ws.Range("K110").Value = 2090.2222
$ws.Range("L110").Value = 1500
$ws.Range("M110").Value = -45.22220000000016
$ws.Range("N110").Value = -5590

# Row 122
$ws.Range("H122").Value = 1591.6538
$ws.Range("I122").Value = 1398.8572
$ws.Range("K122").Value = 4196.571599999999
$ws.Range("M122").Value = -1746.571599999999

# Row 123
$ws.Range("H123").Value = 48429
$ws.Range("J123").Value = 48429
$ws.Range("L123").Value = 48429
$ws.Range("N123").Value = -58229

# Row 136
$ws.Range("H136").Value = 1503.7368
$ws.Range("I136").Value = 1017.25
$ws.Range("J136").Value = 2126.44
$ws.Range("K136").Value = 3051.75
$ws.Range("L136").Value = 6379.32
$ws.Range("M136").Value = -501.75
$ws.Range("N136").Value = -11479.32

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 3802.2
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 3802.2
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 3802.2
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -7296.2

$ws = $wb.Worksheets.Item("CRP")
# Row 42
$ws.Range("H42").Value = 16800
$ws.Range("I42").Value = 10750
$ws.Range("J42").Value = 20833.334
$ws.Range("K42").Value = 10750
$ws.Range("L42").Value = 20833.334
$ws.Range("M42").Value = -10157
$ws.Range("N42").Value = -22019.334

$ws = $wb.Worksheets.Item("CUL")
# Row 117
$ws.Range("H117").Value = 0
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("M117").ClearContents()
$ws.Range("N117").ClearContents()

# Row 120
$ws.Range("H120").Value = 756222.5
$ws.Range("J120").Value = 11000
$ws.Range("L120").Value = 33000
$ws.Range("N120").Value = -42676

# Row 132
$ws.Range("H132").Value = 2242.4443
$ws.Range("I132").Value = 703.6667
$ws.Range("J132").Value = 5320
$ws.Range("K132").Value = 6333.0003
$ws.Range("L132").Value = 47880
$ws.Range("M132").Value = -3803.0003
$ws.Range("N132").Value = -52940

# Row 133
$ws.Range("H133").Value = 6365
$ws.Range("I133").Value = 6365
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 19095
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -14035
$ws.Range("N133").ClearContents()

# Row 134
$ws.Range("H134").Value = 55785812
$ws.Range("I134").Value = 55785812
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 167357436
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -167352366
$ws.Range("N134").ClearContents()

# Row 137
$ws.Range("H137").Value = 13889.889
$ws.Range("I137").Value = 3601.7
$ws.Range("J137").Value = 26750.125
$ws.Range("K137").Value = 10805.1
$ws.Range("L137").Value = 80250.375
$ws.Range("M137").Value = -5705.099999999999
$ws.Range("N137").Value = -90450.375

# Row 139
$ws.Range("H139").Value = 184924.4
$ws.Range("I139").Value = 240900.84
$ws.Range("J139").Value = 3001
$ws.Range("K139").Value = 722702.52
$ws.Range("L139").Value = 9003
$ws.Range("M139").Value = -717562.52
$ws.Range("N139").Value = -19283

# Row 141
$ws.Range("H141").Value = 91185280
$ws.Range("I141").Value = 125377500
$ws.Range("J141").Value = 6000
$ws.Range("K141").Value = 376132500
$ws.Range("L141").Value = 18000
$ws.Range("M141").Value = -376127320
$ws.Range("N141").Value = -28360

$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 11499998
$ws.Range("I11").Value = 11499998
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 11499998
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -11499859
$ws.Range("N11").ClearContents()

# Row 45
$ws.Range("H45").Value = 25000
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 25000
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 25000
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -26118

# Row 70
$ws.Range("H70").Value = 5941.6665
$ws.Range("I70").Value = 5894.737
$ws.Range("J70").Value = 6120
$ws.Range("K70").Value = 5894.737
$ws.Range("L70").Value = 6120
$ws.Range("M70").Value = -5624.737
$ws.Range("N70").Value = -6660

# Row 73
$ws.Range("H73").Value = 5941.6665
$ws.Range("I73").Value = 5894.737
$ws.Range("J73").Value = 6120
$ws.Range("K73").Value = 5894.737
$ws.Range("L73").Value = 6120
$ws.Range("M73").Value = -4958.737
$ws.Range("N73").Value = -7992

# Row 103
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

# Row 126
$ws.Range("H126").Value = 23960
$ws.Range("I126").Value = 133604
$ws.Range("J126").Value = 2031.2
$ws.Range("K126").Value = 400812
$ws.Range("L126").Value = 6093.6
$ws.Range("M126").Value = -398342
$ws.Range("N126").Value = -11033.6

$ws = $wb.Worksheets.Item("LTW")
# Row 127
$ws.Range("H127").Value = 54000
$ws.Range("J127").Value = 54000
$ws.Range("L127").Value = 54000
$ws.Range("N127").Value = -63920

$ws = $wb.Worksheets.Item("WVR")
# Row 123
$ws.Range("H123").Value = 40000
$ws.Range("J123").Value = 40000
$ws.Range("L123").Value = 40000
$ws.Range("N123").Value = -49800

# Row 140
$ws.Range("H140").Value = 41360.832
$ws.Range("J140").Value = 41360.832
$ws.Range("L140").Value = 41360.832
$ws.Range("N140").Value = -51720.832

# Row 141
$ws.Range("H141").Value = 25539.2
$ws.Range("J141").Value = 25539.2
$ws.Range("L141").Value = 25539.2
$ws.Range("N141").Value = -35899.2
